# Update StructureDefinition-mindfulness-audit-retention.xlsx
# - refresh the canonical URL (Metadata!B2 and Elements!R5, which holds the
#   Extension.url fixed value == the StructureDefinition's own canonical URL)
# - refresh the generation Date (Metadata!B8)
# - refresh the "best fit" column widths on the Elements sheet to match the
#   re-rendered IG Publisher output

$wb = $excel.ActiveWorkbook

$newUrl = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/mindfulness-audit-retention"
$newDate = "2025-08-20T10:40:04+01:00"

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = $newUrl
$meta.Range("B8").Value = $newDate

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("R5").Value = $newUrl

# Column widths (Elements sheet) - values are the Excel "ColumnWidth" inputs
# that reproduce the target stored widths (target - 5/6 padding offset).
$widths = @{
    1  = 15.584635416666666
    2  = 15.584635416666666
    3  = 8.959635416666666
    4  = 6.213541666666667
    5  = 4.467447916666667
    6  = 3.1197916666666665
    7  = 3.4322916666666665
    8  = 11.854166666666666
    9  = 9.678385416666666
    11 = 7.565104166666667
    15 = 11.428385416666666
    20 = 6.967447916666667
    21 = 12.776041666666666
    22 = 13.084635416666666
    23 = 14.178385416666666
    24 = 13.795572916666666
    25 = 16.248697916666668
    26 = 14.350260416666666
    27 = 4.240885416666667
    28 = 17.147135416666668
    29 = 33.744791666666664
    30 = 12.709635416666666
    31 = 10.486979166666666
    32 = 14.213541666666666
    33 = 7.389322916666667
    34 = 7.697916666666667
    37 = 18.729166666666668
}

foreach ($col in $widths.Keys) {
    $elements.Columns.Item($col).ColumnWidth = $widths[$col]
}

# Columns 3, 4, 31, 32, 33 are hidden in both before/after states; re-assert
# since setting ColumnWidth above can drop the Hidden flag.
$hiddenCols = @(3, 4, 31, 32, 33)
foreach ($col in $hiddenCols) {
    $elements.Columns.Item($col).Hidden = $true
}
